$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label LF_FFR -> LF_C
$ws.Range("C1").Value = "LF_C"

# Update params row
$ws.Range("B2").Value = -30.22662239492422
$ws.Range("C2").Value = -9.112975023283425

# Update pvalue row
$ws.Range("B3").Value = [double]"1.662652948652976E-08"
$ws.Range("C3").Value = 0.2087492853876425
